$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.176.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.323.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.901.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.319.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.312.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "27.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("E38").Value = "  -6.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.832.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.791"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "322.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0270"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  -2.91%  "
